$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Move the _GoBack bookmark away from "Endre bruk i NOT_MOVING_AT_FLOOR
#    og NOT_MOVING_BETWEEN_FLOORS" and onto the end of the paragraph
#    "Endre plassering av funksjoner ihht tips fra time med Kolbjorn".
# ------------------------------------------------------------------

# Remove the existing (hidden) _GoBack bookmark.
$d.Bookmarks("_GoBack").Delete()

# Locate the "... Kolbjorn" paragraph that should receive the bookmark.
$kolbjornPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*tips fra time med Kolbj*") {
        $kolbjornPara = $p
    }
}

# Placing a zero-length bookmark exactly one character before a paragraph
# mark mis-positions it in this runtime, so temporarily extend the
# paragraph by one placeholder character, drop the bookmark right before
# that placeholder, then remove the placeholder again.
$insertPos = $kolbjornPara.Range.End - 1
$placeholderRange = $d.Range($insertPos, $insertPos)
$placeholderRange.InsertAfter("X")

$bmPos = $insertPos
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$cleanupRange = $d.Range($kolbjornPara.Range.End - 2, $kolbjornPara.Range.End - 1)
$cleanupRange.Delete()

# ------------------------------------------------------------------
# 2) Insert a new to-do item "Legge til timer-modul" right after the
#    "... Kolbjorn" paragraph, matching its list formatting.
# ------------------------------------------------------------------

$kolbjornPara.Range.InsertParagraphAfter()

$newPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*tips fra time med Kolbj*") {
        $newPara = $p.Next(1)
    }
}
$newPara.Range.Text = "Legge til timer-modul"

# ------------------------------------------------------------------
# 3) Highlight the "Fikse ett eller annet rundt case AT_FLOOR i esm.c..."
#    paragraph in yellow.
# ------------------------------------------------------------------

$fiksePara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Fikse ett eller annet rundt*") {
        $fiksePara = $p
    }
}
$fiksePara.Range.HighlightColorIndex = 7
